$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values stay as text (matches original inlineStr formatting)
$ws.Range("D2:D50").NumberFormat = "@"

$ws.Range("D2").Value = '36.994.32'
$ws.Range("D3").Value = '1.979.06'
$ws.Range("E3").Value = '  +1.13%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '244.66'
$ws.Range("E5").Value = '  +0.27%  '
$ws.Range("E6").Value = '  +2.13%  '
$ws.Range("D7").Value = '61.08'
$ws.Range("E7").Value = '  +4.16%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("D10").Value = '0.0798'
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").Value = '14.88'
$ws.Range("E12").Value = '  +8.91%  '
$ws.Range("D13").Value = '22.23'
$ws.Range("E13").Value = '  +0.51%  '
$ws.Range("D14").Value = '0.841'
$ws.Range("E14").Value = '  +1.66%  '
$ws.Range("D15").Value = '2.270.77'
$ws.Range("E15").Value = '  +1.21%  '
$ws.Range("D16").Value = '5.44'
$ws.Range("E16").Value = '  +3.36%  '
$ws.Range("D17").Value = '1.984.09'
$ws.Range("E17").Value = '  +1.10%  '
$ws.Range("D18").Value = '36.878.06'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("D19").Value = '70.06'
$ws.Range("E19").Value = '  +0.58%  '
$ws.Range("D20").Value = '0.0₃0858'
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").Value = '5.14'
$ws.Range("E21").Value = '  +2.02%  '
$ws.Range("D22").Value = '230.15'
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").Value = '2.51'
$ws.Range("E24").Value = '  +2.65%  '
$ws.Range("D25").Value = '2.35'
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '0.149'
$ws.Range("E26").Value = '  +8.71%  '
$ws.Range("D27").Value = '9.26'
$ws.Range("E27").Value = '  +0.24%  '
$ws.Range("D28").Value = '163.02'
$ws.Range("E28").Value = '  +1.74%  '
$ws.Range("D29").Value = '19.59'
$ws.Range("E29").Value = '  +0.95%  '
$ws.Range("E30").Value = '  +16.61%  '
$ws.Range("E31").Value = '  +1.60%  '
$ws.Range("D32").Value = '4.84'
$ws.Range("E32").Value = '  +3.16%  '
$ws.Range("E33").Value = '  +0.40%  '
$ws.Range("D34").Value = '4.54'
$ws.Range("E34").Value = '  +5.80%  '
$ws.Range("E35").Value = '  +2.37%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("D38").Value = '3.33'
$ws.Range("E38").Value = '  -0.23%  '
$ws.Range("D39").Value = '5.48'
$ws.Range("E39").Value = '  -5.03%  '
$ws.Range("D40").Value = '0.0976'
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("E41").Value = '  +1.25%  '
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("D43").Value = '0.0213'
$ws.Range("E43").Value = '  +0.86%  '
$ws.Range("E44").Value = '  +3.81%  '
$ws.Range("D45").Value = '1.368.62'
$ws.Range("E45").Value = '  +0.52%  '
$ws.Range("D46").Value = '89.93'
$ws.Range("E46").Value = '  +2.47%  '
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("D48").Value = '7.21'
$ws.Range("E48").Value = '  +1.36%  '
$ws.Range("E49").Value = '  -0.16%  '
$ws.Range("E50").Value = '  +6.04%  '
